# Update the "Updated symbol list" crypto rates sheet (GitHub Actions run).
#
# Columns D (Price) and E (Volume(1h)) store plain-looking numeric/percent
# text (e.g. "256.40", "0.42%") as literal strings in the workbook, not as
# real numbers. If we just do $ws.Range(...).Value = "256.40", Excel's
# normal type-inference would silently convert it to the number 256.4
# (losing the trailing zero / percent formatting and changing the cell's
# stored type). To reproduce the same text values Excel is first told,
# for the cells we are about to touch, that the column is Text-formatted
# ("@"), so the assignment is kept as a literal string; the number format
# is then restored to General afterwards so the saved cell style matches
# the original (unstyled) cells exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textRange = $ws.Range("D2:E50")
$textRange.NumberFormat = "@"

# --- Row 2 (BNB) ---
$ws.Range("D2").Value = "256.40"
$ws.Range("E2").Value = "0.42%"

# --- Row 3 (OKB) ---
$ws.Range("E3").Value = "-1.16%"

# --- Row 4 (HuobiToken) ---
$ws.Range("D4").Value = "4.618"
$ws.Range("E4").Value = "-10.99%"

# --- Row 5 (Cronos) ---
$ws.Range("D5").Value = "0.05907"
$ws.Range("E5").Value = "0.84%"

# --- Row 6 (KuCoinToken) ---
$ws.Range("D6").Value = "6.638"
$ws.Range("E6").Value = "-1.16%"

# --- Row 7 (MXToken) ---
$ws.Range("D7").Value = "0.8687"
$ws.Range("E7").Value = "0.00%"

# --- Row 8 (FTXToken) ---
$ws.Range("D8").Value = "0.9428"
$ws.Range("E8").Value = "-1.80%"

# --- Row 9: was WazirX -> now One ---
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").Value = "0.01045"
$ws.Range("E9").Value = "-1.11%"

# --- Row 10: was LiechtensteinCryptoassetsExchange -> now WazirX ---
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1407"
$ws.Range("E10").Value = "-0.11%"

# --- Row 11: was MandalaExchangeToken -> now LiechtensteinCryptoassetsExchange ---
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "0.03736"
$ws.Range("E11").Value = "7.96%"

# --- Row 12: was BitrueCoin -> now MandalaExchangeToken ---
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.07085"
$ws.Range("E12").Value = "-1.14%"

# --- Row 13: was BitMartToken -> now BitrueCoin ---
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.03200"
$ws.Range("E13").Value = "-0.12%"

# --- Row 14: was BitForexToken -> now BitMartToken ---
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09264"
$ws.Range("E14").Value = "0.44%"

# --- Row 15: was One -> now BitForexToken ---
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001534"
$ws.Range("E15").Value = "-0.78%"

# --- Row 16 (TigerCash) ---
$ws.Range("D16").Value = "0.006134"
$ws.Range("E16").Value = "2.25%"

# --- Row 17 (LEO) ---
$ws.Range("D17").Value = "3.513"
$ws.Range("E17").Value = "0.39%"

# --- Row 18 (GateToken) ---
$ws.Range("D18").Value = "3.194"

# --- Row 19 (BTSEToken) ---
$ws.Range("E19").Value = "-0.29%"

# --- Row 20 (BitpandaEcosystemToken) ---
$ws.Range("D20").Value = "0.3130"
$ws.Range("E20").Value = "-1.53%"

# --- Row 21 (ProBitToken) ---
$ws.Range("D21").Value = "0.1283"
$ws.Range("E21").Value = "-0.59%"

# --- Row 22 (MCDex) ---
$ws.Range("D22").Value = "3.853"
$ws.Range("E22").Value = "8.91%"

# --- Row 23 (CoinExToken) ---
$ws.Range("D23").Value = "0.04231"
$ws.Range("E23").Value = "1.23%"

# --- Row 24 (ZBToken) ---
$ws.Range("E24").Value = "-1.33%"

# --- Row 25 (BitKan) ---
$ws.Range("D25").Value = "0.001221"
$ws.Range("E25").Value = "-0.03%"

# --- Row 26 (HotbitToken) ---
$ws.Range("D26").Value = "0.004288"
$ws.Range("E26").Value = "-10.47%"

# --- Row 27 (NitroEx) ---
$ws.Range("D27").Value = "0.0001201"
$ws.Range("E27").Value = "0.06%"

# --- Row 28 (UpBots) ---
$ws.Range("D28").Value = "0.0001502"
$ws.Range("E28").Value = "2.46%"

# --- Row 40 (IDEX) ---
$ws.Range("E40").Value = "0.00%"

# --- Row 41 (KickToken) ---
$ws.Range("D41").Value = "0.006193"
$ws.Range("E41").Value = "9.42%"

# --- Row 42 (BKEXToken) ---
$ws.Range("D42").Value = "0.1100"
$ws.Range("E42").Value = "-0.30%"

# --- Row 43 (CEJI) ---
$ws.Range("D43").Value = "0.002246"
$ws.Range("E43").Value = "-2.33%"

# --- Row 44 (LocalTraders) ---
$ws.Range("E44").Value = "16.64%"

# --- Row 45 (CoinLion) ---
$ws.Range("D45").Value = "0.00005502"
$ws.Range("E45").Value = "1.46%"

# --- Row 46 (Kangarootoken) ---
$ws.Range("E46").Value = "0.06%"

# --- Row 47 (CoinbaseStockToken) ---
$ws.Range("D47").Value = "0.08053"
$ws.Range("E47").Value = "-19.46%"

# --- Row 48 (BOLO) ---
$ws.Range("E48").Value = "7.27%"

# --- Row 49 (CryptobidCoin) ---
$ws.Range("D49").Value = "0.00002101"
$ws.Range("E49").Value = "0.06%"

# --- Row 50 (SpecialPowerGold) ---
$ws.Range("D50").Value = "0.0002001"
$ws.Range("E50").Value = "0.06%"

# Restore the original (unstyled/General) look of the touched cells.
$textRange.Style = "Normal"
